# Scenario 1 "Choose a solution and develop a plan to implement it" section
# currently ends with an empty, bulleted (List Paragraph, level-2) paragraph
# that merely carries the document's "_GoBack" bookmark. The plan text is
# added to that paragraph, followed by a blank line and a new (no longer
# bulleted) trailing paragraph that keeps the bookmark - mirroring what Word
# itself does when you type text then hit Enter twice at the end of the very
# last (empty) item of a list.

$d = $word.ActiveDocument

$lastParagraph = $d.Paragraphs.Last

# Build the replacement OOXML: the original paragraph (now carrying the new
# sentences, in two runs) followed by a blank paragraph and a final, empty
# paragraph that keeps the "_GoBack" bookmark. Only the first paragraph
# keeps the List Paragraph style/numbering - the two new ones, just like
# Word produces when you press Enter on an empty last list item, drop the
# list formatting but keep the same direct character formatting.
$newParagraphsXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="41D46342" w14:textId="77777777" w:rsidR="006760C3" w:rsidRPr="00A616C5" w:rsidRDefault="006760C3" w:rsidP="006760C3" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
      <w:b/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>If the man were to transport the cat or the bag of seed first, either choice would leave the bird paired with an unwanted item. Therefore, the bird must be transported first. He must then return for a second item, the choice of which is unimportant. Whether the man returns with the cat or the bag of seed, he will bring the parrot back with him, keeping it safe once again.</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t xml:space="preserve"> He will then switch out the parrot for the third item, and return again for the parrot.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
      <w:b/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica"/>
      <w:b/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
  </w:pPr>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
'@

$lastParagraph.Range.InsertXML($newParagraphsXml)
